$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36; this shifts existing rows 36-166 down to 37-167
$ws.Rows("36").Insert()

# Populate the newly inserted row 36 with the new record
$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(36, 3).Value = "Maule"
$ws.Range("D36").Value = 44690
$ws.Range("D36").NumberFormat = $ws.Range("D37").NumberFormat
$ws.Cells.Item(36, 5).Value = 7
$ws.Cells.Item(36, 6).Value = 100112017
$ws.Cells.Item(36, 7).Value = "Apio"
$ws.Cells.Item(36, 8).Value = "Americana (o)"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 500
$ws.Cells.Item(36, 11).Value = 7000
$ws.Cells.Item(36, 12).Value = 7000
$ws.Cells.Item(36, 13).Value = 7000
$ws.Cells.Item(36, 14).Value = "$/docena de matas"
$ws.Cells.Item(36, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(36, 16).Value = 1167
$ws.Cells.Item(36, 17).Value = 6
$ws.Cells.Item(36, 18).Value = "Hortaliza"
